$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting so numeric-looking values
# such as "1.002" or "218.73" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.144.60'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.657.05'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('D5').Value = '218.73'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '0.5242'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D8').Value = '0.2659'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').Value = '0.06360'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D11').Value = '0.07694'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').Value = '4.611'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').Value = '1.687.75'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').Value = '1.884.40'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').Value = '0.5631'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '0.0₅8200'
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').Value = '65.51'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = '26.143.35'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '4.664'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = '10.53'
$ws.Range('E21').Value = '  +4.13%  '
$ws.Range('D22').Value = '192.63'
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('D23').Value = '5.964'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').Value = '145.30'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').Value = '0.1199'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').Value = '7.270'
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('D28').Value = '15.97'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').Value = '0.05470'
$ws.Range('E30').Value = '  -4.21%  '
$ws.Range('D32').Value = '3.469'
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').Value = '0.9550'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').Value = '2.779'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').Value = '0.5689'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = '5.880'
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').Value = '0.8340'
$ws.Range('E42').Value = '  -1.32%  '
$ws.Range('D43').Value = '1.028.17'
$ws.Range('E43').Value = '  -2.97%  '
$ws.Range('D44').Value = '101.39'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('D45').Value = '1.795.50'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '57.85'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').Value = '8.048'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('D51').Value = '0.05189'
$ws.Range('E51').Value = '  -4.08%  '
